$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Append a new log row at the bottom of the activity log (row 53)
$newRow = 53

$ws.Cells.Item($newRow, 1).Value = 53
$ws.Cells.Item($newRow, 2).Value = "Login"
$ws.Cells.Item($newRow, 3).Value = "User - jiayu logged in."
$ws.Cells.Item($newRow, 4).Value = "09/05/2022 10:25:20 AM"
